$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 38.33049633333334
$ws.Range("H2").Value = 114.991489
$ws.Range("I2").Value = 0.5317874798120843
$ws.Range("J2").Value = 0.5317874798120843
$ws.Range("M2").Value = 20.950428
$ws.Range("N2").Value = 62.851284
$ws.Range("O2").Value = 0.1336707833832197
$ws.Range("P2").Value = 0.1336707833832197
$ws.Range("Q2").Value = 803.040303635764
$ws.Range("R2").Value = 7227.362732721876
$ws.Range("S2").Value = 0.07108444901986943
$ws.Range("T2").Value = 0.07108444901986943

# Row 3
$ws.Range("G3").Value = 38.33049633333334
$ws.Range("H3").Value = 114.991489
$ws.Range("I3").Value = 0.5317874798120843
$ws.Range("J3").Value = 0.5317874798120843
$ws.Range("O3").Value = 0.1853952473714576
$ws.Range("P3").Value = 0.1853952473714576
$ws.Range("Q3").Value = 1113.780079488129
$ws.Range("R3").Value = 10024.02071539316
$ws.Range("S3").Value = 0.0985908713688054
$ws.Range("T3").Value = 0.0985908713688054

# Row 4
$ws.Range("G4").Value = 38.33049633333334
$ws.Range("H4").Value = 114.991489
$ws.Range("I4").Value = 0.5317874798120843
$ws.Range("J4").Value = 0.5317874798120843
$ws.Range("M4").Value = 4.814181666666666
$ws.Range("N4").Value = 14.442545
$ws.Range("O4").Value = 0.03071609967741316
$ws.Range("P4").Value = 0.03071609967741316
$ws.Range("Q4").Value = 184.5299727221672
$ws.Range("R4").Value = 1660.769754499505
$ws.Range("S4").Value = 0.01633443723710832
$ws.Range("T4").Value = 0.01633443723710832

# Row 5
$ws.Range("G5").Value = 38.33049633333334
$ws.Range("H5").Value = 114.991489
$ws.Range("I5").Value = 0.5317874798120843
$ws.Range("J5").Value = 0.5317874798120843
$ws.Range("M5").Value = 95.02938066666667
$ws.Range("N5").Value = 285.088142
$ws.Range("O5").Value = 0.6063194393038427
$ws.Range("P5").Value = 0.6063194393038427
$ws.Range("Q5").Value = 3642.523327202604
$ws.Range("R5").Value = 32782.70994482344
$ws.Range("S5").Value = 0.3224330865884665
$ws.Range("T5").Value = 0.3224330865884665

# Row 6
$ws.Range("G6").Value = 38.33049633333334
$ws.Range("H6").Value = 114.991489
$ws.Range("I6").Value = 0.5317874798120843
$ws.Range("J6").Value = 0.5317874798120843
$ws.Range("M6").Value = 6.880268666666666
$ws.Range("N6").Value = 20.640806
$ws.Range("O6").Value = 0.04389843026406686
$ws.Range("P6").Value = 0.04389843026406686
$ws.Range("Q6").Value = 263.7241129000149
$ws.Range("R6").Value = 2373.517016100134
$ws.Range("S6").Value = 0.02334463559783465
$ws.Range("T6").Value = 0.02334463559783465

# Row 7
$ws.Range("I7").Value = 0.2073226210890634
$ws.Range("J7").Value = 0.2073226210890634
$ws.Range("M7").Value = 20.950428
$ws.Range("N7").Value = 62.851284
$ws.Range("O7").Value = 0.1336707833832197
$ws.Range("P7").Value = 0.1336707833832197
$ws.Range("Q7").Value = 313.073223628272
$ws.Range("R7").Value = 2817.659012654448
$ws.Range("S7").Value = 0.02771297717403752
$ws.Range("T7").Value = 0.02771297717403752

# Row 8
$ws.Range("I8").Value = 0.2073226210890634
$ws.Range("J8").Value = 0.2073226210890634
$ws.Range("O8").Value = 0.1853952473714576
$ws.Range("P8").Value = 0.1853952473714576
$ws.Range("S8").Value = 0.03843662862250588
$ws.Range("T8").Value = 0.03843662862250588

# Row 9
$ws.Range("I9").Value = 0.2073226210890634
$ws.Range("J9").Value = 0.2073226210890634
$ws.Range("M9").Value = 4.814181666666666
$ws.Range("N9").Value = 14.442545
$ws.Range("O9").Value = 0.03071609967741316
$ws.Range("P9").Value = 0.03071609967741316
$ws.Range("Q9").Value = 71.94083927619334
$ws.Range("R9").Value = 647.46755348574
$ws.Range("S9").Value = 0.006368142294754229
$ws.Range("T9").Value = 0.006368142294754229

# Row 10
$ws.Range("I10").Value = 0.2073226210890634
$ws.Range("J10").Value = 0.2073226210890634
$ws.Range("M10").Value = 95.02938066666667
$ws.Range("N10").Value = 285.088142
$ws.Range("O10").Value = 0.6063194393038427
$ws.Range("P10").Value = 0.6063194393038427
$ws.Range("Q10").Value = 1420.073830697469
$ws.Range("R10").Value = 12780.66447627722
$ws.Range("S10").Value = 0.1257037353737239
$ws.Range("T10").Value = 0.1257037353737239

# Row 11
$ws.Range("I11").Value = 0.2073226210890634
$ws.Range("J11").Value = 0.2073226210890634
$ws.Range("M11").Value = 6.880268666666666
$ws.Range("N11").Value = 20.640806
$ws.Range("O11").Value = 0.04389843026406686
$ws.Range("P11").Value = 0.04389843026406686
$ws.Range("Q11").Value = 102.8154599467813
$ws.Range("R11").Value = 925.3391395210319
$ws.Range("S11").Value = 0.009101137624041805
$ws.Range("T11").Value = 0.009101137624041805

# Row 12
$ws.Range("G12").Value = 8.167063666666666
$ws.Range("H12").Value = 24.501191
$ws.Range("I12").Value = 0.1133077476219524
$ws.Range("J12").Value = 0.1133077476219524
$ws.Range("M12").Value = 20.950428
$ws.Range("N12").Value = 62.851284
$ws.Range("O12").Value = 0.1336707833832197
$ws.Range("P12").Value = 0.1336707833832197
$ws.Range("Q12").Value = 171.103479319916
$ws.Range("R12").Value = 1539.931313879244
$ws.Range("S12").Value = 0.01514593538801453
$ws.Range("T12").Value = 0.01514593538801453

# Row 13
$ws.Range("G13").Value = 8.167063666666666
$ws.Range("H13").Value = 24.501191
$ws.Range("I13").Value = 0.1133077476219524
$ws.Range("J13").Value = 0.1133077476219524
$ws.Range("O13").Value = 0.1853952473714576
$ws.Range("P13").Value = 0.1853952473714576
$ws.Range("Q13").Value = 237.312680241351
$ws.Range("R13").Value = 2135.814122172159
$ws.Range("S13").Value = 0.02100671789947456
$ws.Range("T13").Value = 0.02100671789947457

# Row 14
$ws.Range("G14").Value = 8.167063666666666
$ws.Range("H14").Value = 24.501191
$ws.Range("I14").Value = 0.1133077476219524
$ws.Range("J14").Value = 0.1133077476219524
$ws.Range("M14").Value = 4.814181666666666
$ws.Range("N14").Value = 14.442545
$ws.Range("O14").Value = 0.03071609967741316
$ws.Range("P14").Value = 0.03071609967741316
$ws.Range("Q14").Value = 39.31772817456611
$ws.Range("R14").Value = 353.859553571095
$ws.Range("S14").Value = 0.003480372070179065
$ws.Range("T14").Value = 0.003480372070179065

# Row 15
$ws.Range("G15").Value = 8.167063666666666
$ws.Range("H15").Value = 24.501191
$ws.Range("I15").Value = 0.1133077476219524
$ws.Range("J15").Value = 0.1133077476219524
$ws.Range("M15").Value = 95.02938066666667
$ws.Range("N15").Value = 285.088142
$ws.Range("O15").Value = 0.6063194393038427
$ws.Range("P15").Value = 0.6063194393038427
$ws.Range("Q15").Value = 776.111002108569
$ws.Range("R15").Value = 6984.999018977122
$ws.Range("S15").Value = 0.06870069000692351
$ws.Range("T15").Value = 0.06870069000692353

# Row 16
$ws.Range("G16").Value = 8.167063666666666
$ws.Range("H16").Value = 24.501191
$ws.Range("I16").Value = 0.1133077476219524
$ws.Range("J16").Value = 0.1133077476219524
$ws.Range("M16").Value = 6.880268666666666
$ws.Range("N16").Value = 20.640806
$ws.Range("O16").Value = 0.04389843026406686
$ws.Range("P16").Value = 0.04389843026406686
$ws.Range("Q16").Value = 56.19159224443843
$ws.Range("R16").Value = 505.7243301999459
$ws.Range("S16").Value = 0.004974032257360767
$ws.Range("T16").Value = 0.004974032257360767

# Row 17
$ws.Range("G17").Value = 5.834252333333334
$ws.Range("H17").Value = 17.502757
$ws.Range("I17").Value = 0.08094292121735479
$ws.Range("J17").Value = 0.08094292121735479
$ws.Range("M17").Value = 20.950428
$ws.Range("N17").Value = 62.851284
$ws.Range("O17").Value = 0.1336707833832197
$ws.Range("P17").Value = 0.1336707833832197
$ws.Range("Q17").Value = 122.230083443332
$ws.Range("R17").Value = 1100.070750989988
$ws.Range("S17").Value = 0.01081970368845005
$ws.Range("T17").Value = 0.01081970368845005

# Row 18
$ws.Range("G18").Value = 5.834252333333334
$ws.Range("H18").Value = 17.502757
$ws.Range("I18").Value = 0.08094292121735479
$ws.Range("J18").Value = 0.08094292121735479
$ws.Range("O18").Value = 0.1853952473714576
$ws.Range("P18").Value = 0.1853952473714576
$ws.Range("Q18").Value = 169.527521143077
$ws.Range("R18").Value = 1525.747690287693
$ws.Range("S18").Value = 0.0150064329020599
$ws.Range("T18").Value = 0.0150064329020599

# Row 19
$ws.Range("G19").Value = 5.834252333333334
$ws.Range("H19").Value = 17.502757
$ws.Range("I19").Value = 0.08094292121735479
$ws.Range("J19").Value = 0.08094292121735479
$ws.Range("M19").Value = 4.814181666666666
$ws.Range("N19").Value = 14.442545
$ws.Range("O19").Value = 0.03071609967741316
$ws.Range("P19").Value = 0.03071609967741316
$ws.Range("Q19").Value = 28.08715062184056
$ws.Range("R19").Value = 252.784355596565
$ws.Range("S19").Value = 0.00248625083629327
$ws.Range("T19").Value = 0.00248625083629327

# Row 20
$ws.Range("G20").Value = 5.834252333333334
$ws.Range("H20").Value = 17.502757
$ws.Range("I20").Value = 0.08094292121735479
$ws.Range("J20").Value = 0.08094292121735479
$ws.Range("M20").Value = 95.02938066666667
$ws.Range("N20").Value = 285.088142
$ws.Range("O20").Value = 0.6063194393038427
$ws.Range("P20").Value = 0.6063194393038427
$ws.Range("Q20").Value = 554.4253858897217
$ws.Range("R20").Value = 4989.828473007495
$ws.Range("S20").Value = 0.04907726660812167
$ws.Range("T20").Value = 0.04907726660812167

# Row 21
$ws.Range("G21").Value = 5.834252333333334
$ws.Range("H21").Value = 17.502757
$ws.Range("I21").Value = 0.08094292121735479
$ws.Range("J21").Value = 0.08094292121735479
$ws.Range("M21").Value = 6.880268666666666
$ws.Range("N21").Value = 20.640806
$ws.Range("O21").Value = 0.04389843026406686
$ws.Range("P21").Value = 0.04389843026406686
$ws.Range("Q21").Value = 40.14122352246022
$ws.Range("R21").Value = 361.271011702142
$ws.Range("S21").Value = 0.003553267182429907
$ws.Range("T21").Value = 0.003553267182429907

# Row 22
$ws.Range("G22").Value = 4.803262333333334
$ws.Range("H22").Value = 14.409787
$ws.Range("I22").Value = 0.06663923025954499
$ws.Range("J22").Value = 0.066639230259545
$ws.Range("M22").Value = 20.950428
$ws.Range("N22").Value = 62.851284
$ws.Range("O22").Value = 0.1336707833832197
$ws.Range("P22").Value = 0.1336707833832197
$ws.Range("Q22").Value = 100.630401679612
$ws.Range("R22").Value = 905.6736151165081
$ws.Range("S22").Value = 0.008907718112848137
$ws.Range("T22").Value = 0.008907718112848139

# Row 23
$ws.Range("G23").Value = 4.803262333333334
$ws.Range("H23").Value = 14.409787
$ws.Range("I23").Value = 0.06663923025954499
$ws.Range("J23").Value = 0.066639230259545
$ws.Range("O23").Value = 0.1853952473714576
$ws.Range("P23").Value = 0.1853952473714576
$ws.Range("Q23").Value = 139.569752942907
$ws.Range("R23").Value = 1256.127776486163
$ws.Range("S23").Value = 0.01235459657861187
$ws.Range("T23").Value = 0.01235459657861187

# Row 24
$ws.Range("G24").Value = 4.803262333333334
$ws.Range("H24").Value = 14.409787
$ws.Range("I24").Value = 0.06663923025954499
$ws.Range("J24").Value = 0.066639230259545
$ws.Range("M24").Value = 4.814181666666666
$ws.Range("N24").Value = 14.442545
$ws.Range("O24").Value = 0.03071609967741316
$ws.Range("P24").Value = 0.03071609967741316
$ws.Range("Q24").Value = 23.12377746532389
$ws.Range("R24").Value = 208.113997187915
$ws.Range("S24").Value = 0.002046897239078271
$ws.Range("T24").Value = 0.002046897239078271

# Row 25
$ws.Range("G25").Value = 4.803262333333334
$ws.Range("H25").Value = 14.409787
$ws.Range("I25").Value = 0.06663923025954499
$ws.Range("J25").Value = 0.066639230259545
$ws.Range("M25").Value = 95.02938066666667
$ws.Range("N25").Value = 285.088142
$ws.Range("O25").Value = 0.6063194393038427
$ws.Range("P25").Value = 0.6063194393038427
$ws.Range("Q25").Value = 456.4510447161949
$ws.Range("R25").Value = 4108.059402445754
$ws.Range("S25").Value = 0.04040466072660698
$ws.Range("T25").Value = 0.04040466072660699

# Row 26
$ws.Range("G26").Value = 4.803262333333334
$ws.Range("H26").Value = 14.409787
$ws.Range("I26").Value = 0.06663923025954499
$ws.Range("J26").Value = 0.066639230259545
$ws.Range("M26").Value = 6.880268666666666
$ws.Range("N26").Value = 20.640806
$ws.Range("O26").Value = 0.04389843026406686
$ws.Range("P26").Value = 0.04389843026406686
$ws.Range("Q26").Value = 33.04773532981356
$ws.Range("R26").Value = 297.429617968322
$ws.Range("S26").Value = 0.00292535760239973
$ws.Range("T26").Value = 0.002925357602399731
